$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "3.2 Coding Style"
$ws.Range("A24").Value = "3.3 Comments"
$ws.Range("A25").Value = "3.4 Ninja code"

$ws.Range("A25").Select()
